# Applies the diff described in the commit "test elec sec US files":
#  - About sheet: remove the "Nevada" label and date stamp that had been
#    added in B1/C1, leaving just the title in A1.
#  - LFHVM sheet: rename the first two peak rows from "Summer Peak 1" /
#    "Summer Peak 2" to the generic "Summer Peak" / "Winter Peak", and
#    delete the remaining "Summer/Winter Peak 3-5" rows (rows 8-15),
#    shrinking the table back down to a single Summer Peak/Winter Peak row.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsLfhvm = $wb.Worksheets.Item("LFHVM")

# --- About sheet: drop the state label + timestamp cells in row 1 ---
$wsAbout.Range("B1:C1").Clear()

# --- LFHVM sheet: collapse the 10 peak rows down to 2 ---
$wsLfhvm.Range("A6").Value = "Summer Peak"
$wsLfhvm.Range("A7").Value = "Winter Peak"
$wsLfhvm.Range("A8:I15").Delete()

# Restore the selection/active-sheet state: leave a lingering selection on
# F25 in LFHVM (from editing it) but land back on the About tab.
$null = $wsLfhvm.Range("F25").Select()
$null = $wsAbout.Select()
